$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 317010
$ws.Range("D2").Value = 404006357
$ws.Range("C10").Value = 116209
$ws.Range("D10").Value = 170283149
$ws.Range("C12").Value = 58729
$ws.Range("D12").Value = 84755241
$ws.Range("C16").Value = 3984
$ws.Range("D16").Value = 5653373
$ws.Range("C20").Value = 6524
$ws.Range("D20").Value = 9099822
$ws.Range("C22").Value = 76593
$ws.Range("D22").Value = 95548273
$ws.Range("C28").Value = 32245
$ws.Range("D28").Value = 47207204
$ws.Range("C30").Value = 11377
$ws.Range("D30").Value = 16363039
$ws.Range("C35").Value = 1791
$ws.Range("D35").Value = 2528115
$ws.Range("C36").Value = 96223
$ws.Range("D36").Value = 121158445
$ws.Range("C44").Value = 44105
$ws.Range("D44").Value = 64634677
$ws.Range("C46").Value = 9050
$ws.Range("D46").Value = 12987669
$ws.Range("C51").Value = 2263
$ws.Range("D51").Value = 3157894
$ws.Range("C52").Value = 68341
$ws.Range("D52").Value = 85750065
$ws.Range("C58").Value = 27945
$ws.Range("D58").Value = 40983109
$ws.Range("C61").Value = 10977
$ws.Range("D61").Value = 15872241
$ws.Range("C67").Value = 1446
$ws.Range("D67").Value = 2025065
$ws.Range("C69").Value = 20285
$ws.Range("D69").Value = 26569872
$ws.Range("C73").Value = 7529
$ws.Range("D73").Value = 11023021
$ws.Range("C75").Value = 5065
$ws.Range("D75").Value = 7353678
$ws.Range("C76").Value = 487
$ws.Range("D76").Value = 689239
$ws.Range("C78").Value = 139241
$ws.Range("D78").Value = 173652645
$ws.Range("C84").Value = 63074
$ws.Range("D84").Value = 92445784
$ws.Range("C87").Value = 29415
$ws.Range("D87").Value = 42547516
$ws.Range("C89").Value = 2721
$ws.Range("D89").Value = 3917852
$ws.Range("C91").Value = 32399
$ws.Range("D91").Value = 43887633
$ws.Range("C95").Value = 7836
$ws.Range("D95").Value = 11521581
$ws.Range("C97").Value = 7155
$ws.Range("D97").Value = 10373082
$ws.Range("C99").Value = 527
$ws.Range("D99").Value = 749905
$ws.Range("C100").Value = 483
$ws.Range("D100").Value = 696943
$ws.Range("C101").Value = 8921
$ws.Range("D101").Value = 12384334
$ws.Range("C103").Value = 2240
$ws.Range("D103").Value = 3300302
$ws.Range("C105").Value = 3010
$ws.Range("D105").Value = 4396550
$ws.Range("C107").Value = 131
$ws.Range("D107").Value = 190620
$ws.Range("C109").Value = 139696
$ws.Range("D109").Value = 172743359
$ws.Range("C115").Value = 52298
$ws.Range("D115").Value = 76665589
$ws.Range("C117").Value = 26677
$ws.Range("D117").Value = 38647996
$ws.Range("C118").Value = 1304
$ws.Range("D118").Value = 1784458
$ws.Range("C121").Value = 2211
$ws.Range("D121").Value = 3105355
$ws.Range("C123").Value = 495760
$ws.Range("D123").Value = 653755340
$ws.Range("C125").Value = 209
$ws.Range("D125").Value = 308236
$ws.Range("C130").Value = 205380
$ws.Range("D130").Value = 301906259
$ws.Range("C133").Value = 177596
$ws.Range("D133").Value = 258135446
$ws.Range("C136").Value = 2833
$ws.Range("D136").Value = 3982381
$ws.Range("C138").Value = 6208
$ws.Range("D138").Value = 8771649
$ws.Range("C141").Value = 43975
$ws.Range("D141").Value = 58702970
$ws.Range("C147").Value = 13933
$ws.Range("D147").Value = 20432676
$ws.Range("C148").Value = 3710
$ws.Range("D148").Value = 5349680
$ws.Range("C151").Value = 397
$ws.Range("D151").Value = 571431
$ws.Range("C154").Value = 17344
$ws.Range("D154").Value = 22914662
$ws.Range("C158").Value = 7085
$ws.Range("D158").Value = 10304703
$ws.Range("C160").Value = 4931
$ws.Range("D160").Value = 7095563
$ws.Range("C162").Value = 275
$ws.Range("D162").Value = 380231
$ws.Range("C163").Value = 264
$ws.Range("D163").Value = 377364
$ws.Range("C165").Value = 15649
$ws.Range("D165").Value = 22705655
$ws.Range("C166").Value = 1760
$ws.Range("D166").Value = 2617730
$ws.Range("C167").Value = 236
$ws.Range("D167").Value = 348802
$ws.Range("C171").Value = 86780
$ws.Range("D171").Value = 108550899
$ws.Range("C178").Value = 33598
$ws.Range("D178").Value = 49270321
$ws.Range("C180").Value = 12865
$ws.Range("D180").Value = 18586409
$ws.Range("C182").Value = 1241
$ws.Range("D182").Value = 1737396
$ws.Range("C184").Value = 1618
$ws.Range("D184").Value = 2272662
$ws.Range("C186").Value = 235963
$ws.Range("D186").Value = 293325836
$ws.Range("C194").Value = 85993
$ws.Range("D194").Value = 126053150
$ws.Range("C197").Value = 32708
$ws.Range("D197").Value = 47072136
$ws.Range("C200").Value = 5084
$ws.Range("D200").Value = 7241993
$ws.Range("C203").Value = 4787
$ws.Range("D203").Value = 6625058
$ws.Range("C206").Value = 260932
$ws.Range("D206").Value = 322932991
$ws.Range("C215").Value = 94440
$ws.Range("D215").Value = 138157686
$ws.Range("C218").Value = 50906
$ws.Range("D218").Value = 73567874
$ws.Range("C221").Value = 4651
$ws.Range("D221").Value = 6529876
$ws.Range("C224").Value = 5641
$ws.Range("D224").Value = 7798448
$ws.Range("C227").Value = 105043
$ws.Range("D227").Value = 131425378
$ws.Range("C234").Value = 49134
$ws.Range("D234").Value = 71980905
$ws.Range("C236").Value = 12244
$ws.Range("D236").Value = 17603708
$ws.Range("C240").Value = 2457
$ws.Range("D240").Value = 3433565
$ws.Range("C241").Value = 254421
$ws.Range("D241").Value = 321218721
$ws.Range("C242").Value = 172
$ws.Range("D242").Value = 212933
$ws.Range("C249").Value = 94998
$ws.Range("D249").Value = 139194002
$ws.Range("C252").Value = 64196
$ws.Range("D252").Value = 93024138
$ws.Range("C254").Value = 2396
$ws.Range("D254").Value = 3380928
$ws.Range("C257").Value = 4519
$ws.Range("D257").Value = 6344116
